$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.210.48"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.267.78"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.493"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.45%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "2.619.41"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "2.276.45"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.781"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "42.139.22"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0683"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").Value = "1.968.02"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0276"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("E48").Value = "  -4.80%  "
$ws.Range("D49").Value = "2.492.68"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("E51").Value = "  -0.56%  "
